# Add a new submission row (row 28) to the first worksheet
# ("八位序列号收集收集结果yd5") of the "八位序列号收集（收集结果）" workbook.
#
# Columns: A=提交者（自动） B=提交时间（自动） C=序列号（必填） D=QQ号（必填）

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# A28 - submitter name (plain text)
$ws.Range("A28").Value = "二十一"

# B28 - submission timestamp, stored as a date serial with the same
# "yyyy/m/d h:mm:ss;@" format used by the rows above it.
$ws.Range("B28").Value = 45902.3671412037
$ws.Range("B28").NumberFormat = "yyyy/m/d h:mm:ss;@"

# C28 - serial number (plain text)
$ws.Range("C28").Value = "e964b1e9"

# D28 - QQ number. It is all digits, so force it to be stored as text
# (matching the rest of column D) instead of being auto-coerced to a
# number, then drop the temporary text format so the cell keeps the
# workbook's default (unstyled) look.
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "820078578"
$ws.Range("D28").ClearFormats()
